$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "70.378.97"
$ws.Range("E2").Value = "  +4.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.627.45"
$ws.Range("E3").Value = "  +4.45%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "591.35"
$ws.Range("E5").Value = "  +0.98%  "

# Row 6 - Solana
$ws.Range("D6").Value = "194.22"
$ws.Range("E6").Value = "  +3.97%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +2.00%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.621.45"
$ws.Range("E8").Value = "  +4.34%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  -0.02%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +4.78%  "

# Row 11 - Cardano
$ws.Range("D11").Value = "0.671"
$ws.Range("E11").Value = "  +3.26%  "

# Row 12 - Avalanche
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.40"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.85%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +4.30%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "9.96"
$ws.Range("E14").Value = "  +5.63%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.209.67"
$ws.Range("E15").Value = "  +4.32%  "

# Row 16 - Chainlink
$ws.Range("D16").Value = "19.88"
$ws.Range("E16").Value = "  +5.59%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.625.70"
$ws.Range("E17").Value = "  +4.44%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "70.335.94"
$ws.Range("E18").Value = "  +4.68%  "

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +4.52%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +1.96%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  +4.37%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "488.05"
$ws.Range("E22").Value = "  +0.46%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("D23").Value = "18.97"
$ws.Range("E23").Value = "  +12.53%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "5.31"

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "4.46"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26 - Litecoin
$ws.Range("D26").Value = "91.61"
$ws.Range("E26").Value = "  +1.94%  "

# Row 27 - ImmutableX
$ws.Range("D27").Value = "3.17"
$ws.Range("E27").Value = "  +7.30%  "

# Row 28 - RenderToken
$ws.Range("E28").Value = "  +5.09%  "

# Row 29 - Filecoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.60"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.94%  "

# Row 30 - EthereumClassic
$ws.Range("D30").Value = "32.95"
$ws.Range("E30").Value = "  +5.05%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  +9.45%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  +8.29%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "626.93"
$ws.Range("E33").Value = "  +4.14%  "

# Row 34 - Cosmos
$ws.Range("E34").Value = "  +4.45%  "

# Row 35 - OKB
$ws.Range("D35").Value = "65.63"
$ws.Range("E35").Value = "  +2.90%  "

# Row 36 - InjectiveProtocol
$ws.Range("D36").Value = "40.79"
$ws.Range("E36").Value = "  +11.50%  "

# Row 37 - TheGraph
$ws.Range("E37").Value = "  +6.79%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0823"
$ws.Range("E38").Value = "  +8.32%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  -1.05%  "

# Row 40 - Dai
$ws.Range("E40").Value = "  +0.00%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +1.17%  "

# Row 42 - Maker
$ws.Range("D42").Value = "3.299.34"
$ws.Range("E42").Value = "  +0.89%  "

# Row 43 - was ThetaToken, now Fetch.AI
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D43").Value = "2.89"
$ws.Range("E43").Value = "  +14.20%  "

# Row 44 - was Fetch.AI, now ThetaToken
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").Value = "3.19"
$ws.Range("E44").Value = "  +9.53%  "

# Row 45 - VeChain
$ws.Range("E45").Value = "  +5.58%  "

# Row 46 - dogwifhat
$ws.Range("D46").Value = "2.91"
$ws.Range("E46").Value = "  +3.10%  "

# Row 47 - ApeXProtocol
$ws.Range("E47").Value = "  +1.05%  "

# Row 48 - Stellar
$ws.Range("E48").Value = "  +2.37%  "

# Row 49 - THORChain
$ws.Range("D49").Value = "9.21"
$ws.Range("E49").Value = "  +5.21%  "

# Row 50 - LidoDAOToken
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  +1.71%  "

# Row 51 - FirstDigitalUSD
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.00%  "
